# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which data-row (rows 2-20) each market
# record lives on. Columns A,B,C,E,F,G,H,I,J,K,R are constant across all
# rows already, so only D (Fecha), L (Calidad), M (Volumen), N (Precio
# minimo), O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), S (Precio $/Kg) and T (Kg / unidad) actually need to
# move between rows.
#
# Mapping below: for each destination row (key), the value is the row
# that currently (pre-edit) holds the data that should land there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 5
    3  = 8
    4  = 18
    5  = 12
    6  = 15
    7  = 16
    8  = 17
    9  = 20
    10 = 6
    11 = 9
    12 = 2
    13 = 3
    14 = 4
    15 = 19
    16 = 13
    17 = 14
    18 = 10
    19 = 11
    20 = 7
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot every value we might need to move, BEFORE any writes, so
# overlapping cycles in the permutation don't clobber source data.
$snapshot = @{}
for ($row = 2; $row -le 20; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $snapshot["$col$srcRow"]
    }
}
